# Generate Report for Handoff
# Updates the localization-status workbook after a new handoff report run:
#  - Priority column ("low" -> "ht") for the md/markdown rows that were handed off
#  - Latest Handoff Datetime column is refreshed with the new handoff timestamps

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4:G7").Value = "2016-08-25 16:33:35"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4").Value = "2016-08-25 16:33:30"
$zhcn.Range("H5").Value = "2016-08-25 16:33:30"
$zhcn.Range("H6").Value = "2016-08-25 16:33:30"
$zhcn.Range("H7").Value = "2016-08-25 16:33:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4").Value = "2016-08-25 16:33:35"
$dede.Range("H5").Value = "2016-08-25 16:33:35"
$dede.Range("H6").Value = "2016-08-25 16:33:35"
$dede.Range("H7").Value = "2016-08-25 16:33:35"
